# docs: update calculator_documentation.docx for feat: added function - added factorial

$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. Table of Contents entry "4. Feature Specifications" gains a new
#    bullet line (via a manual line break) describing the new
#    Factorial Operation feature.
# -----------------------------------------------------------------
$tocRange = $d.Content
$found = $tocRange.Find.Execute("4. Feature Specifications", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $tocRange.Collapse(0)
    $bullet = [char]8226
    $newLine = $bullet + " Factorial Operation: Calculates the factorial of a number"
    $tocRange.InsertAfter([char]11 + $newLine)
}

# -----------------------------------------------------------------
# 2. Menu options table (Option | Function) gains a new row
#    describing the Exponentiation / power function, appended right
#    after the existing "Exit" row.
# -----------------------------------------------------------------
$table = $d.Tables.Item(1)
$newRow = $table.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "Exponentiation | power(a, b)"
$newRow.Cells.Item(2).Range.Text = "power"
